# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect freshly generated output, per commit "Update gh-pages to
# output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Row => new value for column F, identical update set for both sheets
# (sheet "全部类型" has one extra row inserted around row 20, so the
# row numbers differ slightly between the two sheets).
$changesExhibition = @{
    2  = 15141
    3  = 19423
    5  = 159
    14 = 202
    15 = 241
    16 = 72
    17 = 1504
    20 = 108
    21 = 244
    22 = 8140
    23 = 991
    25 = 10
    27 = 1268
    28 = 4
    31 = 6515
    35 = 157
    36 = 298
    37 = 5539
    38 = 1014
    39 = 25
    41 = 59
}

$changesAllTypes = @{
    2  = 15141
    3  = 19423
    5  = 159
    14 = 202
    15 = 241
    16 = 72
    17 = 1504
    21 = 108
    22 = 244
    23 = 8140
    24 = 991
    26 = 10
    28 = 1268
    29 = 4
    34 = 6515
    38 = 157
    39 = 298
    40 = 5539
    41 = 1014
    42 = 25
    44 = 59
}

$ws = $wb.Worksheets.Item("展览")
foreach ($row in $changesExhibition.Keys) {
    $ws.Range("F$row").Value = $changesExhibition[$row]
}

$ws = $wb.Worksheets.Item("全部类型")
foreach ($row in $changesAllTypes.Keys) {
    $ws.Range("F$row").Value = $changesAllTypes[$row]
}
